$d = $word.ActiveDocument

# Event date
$d.Content.Find.Execute("04/septiembre/2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "13/septiembre/2017", 2)

# Consumer name (appears twice in the document)
$d.Content.Find.Execute("Jose Perez Martinez", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ofelia  Martinez Zamora", 2)

# Event day
$d.Content.Find.Execute("07/septiembre/2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "15/septiembre/2017", 2)

# Time
$d.Content.Find.Execute("13:00hrs.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "12:00hrs.", 2)

# Number of people (whole word match to avoid collateral hits)
$d.Content.Find.Execute("19", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "150", 2)

# Place
$d.Content.Find.Execute("Casa de Pancha", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "En un lugar bien chido", 2)

# Dish
$d.Content.Find.Execute("Pollo con Piña, 3 tiempos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Platillo de prueba, 3 tiempos", 2)

# Services included
$d.Content.Find.Execute("4 Mesero, 2 Platos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2 Platos, 1 Cubiertos, 1 Vasos High Ball, 4 Barra de Licores y Cerveza", 2)

# Price
$d.Content.Find.Execute("`$13,000.00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "`$30,000.00", 2)

# Deposit
$d.Content.Find.Execute("`$2,000.00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "`$1,500.00", 2)

# Notes
$d.Content.Find.Execute("Comida de cconvivio empresarial", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Una descripcion muy completa", 2)
